$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a text code ("002" -> "004"). Force text so Excel doesn't coerce it
# to a number (which would drop the leading zeros), then strip the
# "@" number-format style back off so the cell keeps its original
# (unstyled) appearance.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"
$ws.Range("J2").Style = "Normal"

# N2 is a plain text timestamp string.
$ws.Range("N2").Value = "2017-09-30 00:00:00"

# Numeric cells.
$ws.Range("O2").Value = -299059278.04
$ws.Range("P2").Value = -294.1697714549
$ws.Range("Q2").Value = 505508660.84
$ws.Range("R2").Value = 497.2437845847
$ws.Range("S2").Value = 42286692.73
$ws.Range("T2").Value = 41.5953212269
$ws.Range("U2").Value = -2858588.48
$ws.Range("V2").Value = -2.8118516348
$ws.Range("Y2").Value = 2858588.48
$ws.Range("Z2").Value = 2.8118516348
$ws.Range("AA2").Value = 200217651.84
$ws.Range("AB2").Value = 196.9441686244
$ws.Range("AC2").Value = -101662137.67
$ws.Range("AD2").Value = -100.6589334423
